# Tidy up label for rmd: rename "Urbanicity" label to "Urbanicity (Rur.)"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Update the "label" value for the urban_rura_fctb row (was "Urbanicity")
$ws.Range("B13").Value = "Urbanicity (Rur.)"

# Move the active selection to B14, matching the end state after the edit
$ws.Range("B14").Select()
